# Optimize database schema migration
# Append one new data row to each of the four worksheets, mirroring the
# existing row layout (time / 总长 / ID / 实际长度 / 和校验 / ..._DEC columns).
# NOTE: this engine's PowerShell subset does not reliably bind "-Name value"
# style named arguments, and does not like inline parenthesised expressions
# as bare call arguments either, so every argument is first materialised
# into its own variable and Add-DataRow is always called positionally.

function Add-DataRow {
    param($ws, $row, $timeValue, $colB, $colC, $colD, $colE, $colF, $colG, $colH, $colI)

    # Column A: timestamp - reuse the date/time number format already used
    # by the row directly above so the new cell matches the existing style.
    $prevRow = $row - 1
    $aCell = $ws.Cells.Item($row, 1)
    $aFmt = $ws.Cells.Item($prevRow, 1).NumberFormat
    $aCell.NumberFormat = $aFmt
    $aCell.Value = $timeValue

    # Columns B-E: free-form hex strings, stored as plain text.
    $ws.Cells.Item($row, 2).Value = $colB
    $ws.Cells.Item($row, 3).Value = $colC
    $ws.Cells.Item($row, 4).Value = $colD
    $ws.Cells.Item($row, 5).Value = $colE

    # Column F: numeric length value.
    $ws.Cells.Item($row, 6).Value = $colF

    # Column G: decimal ID value. Some rows store it as a genuine number;
    # others (when it overflows reliable double precision) keep it as text
    # so no precision is lost - the caller passes a [double] or a [string]
    # accordingly.
    $gCell = $ws.Cells.Item($row, 7)
    $gIsString = $colG -is [string]
    if ($gIsString) {
        $gCell.NumberFormat = "@"
        $gCell.Value = $colG
        # Drop the temporary "Text" number format now that the value is
        # safely stored as a string, so the cell is left with the default
        # (unstyled) appearance instead of permanently showing "@".
        $gCell.Style = "Normal"
    } else {
        $gCell.Value = $colG
    }

    # Columns H-I: numeric values.
    $ws.Cells.Item($row, 8).Value = $colH
    $ws.Cells.Item($row, 9).Value = $colI
}

$wb = $excel.ActiveWorkbook

# --- Sheet 1: ROW50-FE-LIFTER -> new row 48 ---
$ws1 = $wb.Worksheets.Item("ROW50-FE-LIFTER")
$row1 = 48
$time1 = [double]"45748.18585129629"
$b1 = "0x01,0x90"
$c1 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$d1 = "0x01,0x66"
$e1 = "0xe"
$f1 = 400
$g1 = [double]"5.68631262647114e+23"
$h1 = 358
$i1 = 14
Add-DataRow $ws1 $row1 $time1 $b1 $c1 $d1 $e1 $f1 $g1 $h1 $i1

# --- Sheet 2: ROW50-MID-LIFTER -> new row 50 ---
$ws2 = $wb.Worksheets.Item("ROW50-MID-LIFTER")
$row2 = 50
$time2 = [double]"45748.15712962963"
$b2 = "0x01,0x90 "
$c2 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$d2 = "0x01,0x6a"
$e2 = "0x19"
$f2 = 400
$g2 = "568631262647113771663628"
$h2 = 362
$i2 = 25
Add-DataRow $ws2 $row2 $time2 $b2 $c2 $d2 $e2 $f2 $g2 $h2 $i2

# --- Sheet 3: ROW11-FE-LIFTER -> new row 48 ---
$ws3 = $wb.Worksheets.Item("ROW11-FE-LIFTER")
$row3 = 48
$time3 = [double]"45748.21298219907"
$b3 = "0x01,0x90"
$c3 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$d3 = "0x01,0x66"
$e3 = "0x14"
$f3 = 400
$g3 = [double]"5.68631262647114e+23"
$h3 = 358
$i3 = 20
Add-DataRow $ws3 $row3 $time3 $b3 $c3 $d3 $e3 $f3 $g3 $h3 $i3

# --- Sheet 4: ROW11-MID-LIFTER -> new row 48 ---
$ws4 = $wb.Worksheets.Item("ROW11-MID-LIFTER")
$row4 = 48
$time4 = [double]"45748.34986975694"
$b4 = "0x01,0x90"
$c4 = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$d4 = "0x01,0x6e"
$e4 = "0x19"
$f4 = 400
$g4 = [double]"5.68631262647114e+23"
$h4 = 366
$i4 = 25
Add-DataRow $ws4 $row4 $time4 $b4 $c4 $d4 $e4 $f4 $g4 $h4 $i4

Write-Host "Appended new rows to all four worksheets."
